# Mifos -> Finflux: 1st changes
# - Insert a new blank column before the "Late" column on the "Repayment schedule" sheet
#   (shifts Late / Paid Date / Outstanding one column to the right).
# - Make "Repayment schedule" the active sheet/tab (was "Transactions").
# - Update the selection on "Repayment schedule" to S5.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at column N (14), pushing the existing
# "Late" (N), "Paid Date" (O) and "Outstanding" (P) columns right by one.
$wsSchedule.Range("N1").EntireColumn.Insert() | Out-Null

# Excel copies the format (incl. width) of the column to the left when inserting;
# reproduce the resulting column width as closely as possible.
$wsSchedule.Range("N1").ColumnWidth = $wsSchedule.Range("M1").ColumnWidth

# Make the "Repayment schedule" sheet the active tab (previously "Transactions" was active).
$wsSchedule.Activate() | Out-Null

# Update the active selection on the now-active sheet.
$wsSchedule.Range("S5").Select() | Out-Null
